# Insert a new weekly price record as row 63 in the Ají (Chillán) sheet.
# This pushes the existing rows 63-127 down to 64-128 (Excel's native
# Rows.Insert shifts formatting + values together), then the new row is
# populated with the latest "Americana (o)" observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(63).Insert()

$ws.Range("A63").Value = 7
$ws.Range("B63").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C63").Value = "Ñuble"
$ws.Range("D63").Value = 44944
$ws.Range("E63").Value = 16
$ws.Range("F63").Value = 100112021
$ws.Range("G63").Value = "Ají"
$ws.Range("H63").Value = "Americana (o)"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 50
$ws.Range("K63").Value = 13000
$ws.Range("L63").Value = 13000
$ws.Range("M63").Value = 13000
$ws.Range("N63").Value = "`$/caja 15 kilos"
$ws.Range("O63").Value = "Región del Maule"
$ws.Range("P63").Value = 867
$ws.Range("Q63").Value = 15
$ws.Range("R63").Value = "Hortaliza"
